$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data in columns D and E (header text + all observation values),
# row by row, for rows 1 (header) through 67 (last data row).
for ($r = 1; $r -le 67; $r++) {
    $d = $ws.Cells.Item($r, 4).Value2
    $e = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 4).Value2 = $e
    $ws.Cells.Item($r, 5).Value2 = $d
}

# Column widths for D and E changed alongside the data swap.
$ws.Columns.Item(4).ColumnWidth = 14.497395833333332
$ws.Columns.Item(5).ColumnWidth = 15.997395833333332

# Update the view: zoom level and active selection.
$excel.ActiveWindow.Zoom = 62
$null = $ws.Range("L13").Select()
